$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" ---
# This string appears on the Overview sheet (language status columns) and
# on each per-language detail sheet ("Status" column), so every occurrence
# needs to be updated.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Column width changes (narrower "Status"/status-date columns) ---
# Overview sheet: columns E and F (zh-cn / de-de status columns)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C ("Status")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C ("Status")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
